$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2 (first data row): time + comprador
$t.Cell(2, 1).Range.Text = "08:45 - 09:00"
$t.Cell(2, 3).Range.Text = "BOX BRAND"

# Row 3
$t.Cell(3, 1).Range.Text = "09:00 - 09:15"
$t.Cell(3, 3).Range.Text = "INTERLINK2AMERICAS"

# Row 4
$t.Cell(4, 3).Range.Text = "INMERSSO BOUTIQUE"

# Row 5
$t.Cell(5, 3).Range.Text = "COLFRESH COFFEE"

# Row 6
$t.Cell(6, 1).Range.Text = "11:00 - 11:15"
